$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FBHS")

# Row 4 - Inventory
$ws.Range("C4").Value = 867000000.0
$ws.Range("D4").Value = 738000000.0
$ws.Range("E4").Value = 732000000.0
$ws.Range("F4").Value = 703000000.0
$ws.Range("G4").Value = 719000000.0

# Row 14 - Accounts Payable
$ws.Range("C14").Value = 621000000.0
$ws.Range("D14").Value = 544000000.0
$ws.Range("E14").Value = 464000000.0
$ws.Range("F14").Value = 426000000.0
$ws.Range("G14").Value = 460000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("C24").Value = 130000000.0
$ws.Range("D24").Value = 149000000.0
$ws.Range("E24").Value = 152000000.0
$ws.Range("F24").Value = 157000000.0
$ws.Range("G24").Value = 140000000.0

# Row 40 - Net Debt (B40 goes from empty inline string to a numeric value)
$ws.Range("B40").Value = 2326700000.0

# Row 41 - Total Debt (B41 goes from empty inline string to a numeric value)
$ws.Range("B41").Value = 2682800000.0
